$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data for "Reverse a Linked List"
$ws.Range("A8").Value = "Linked List"
$ws.Range("B8").Value = "Reverse a Linked List"
$ws.Range("C8").Value = "https://leetcode.com/explore/interview/card/top-interview-questions-easy/93/linked-list/560/"
$ws.Range("F8").Value = "Easy - keep track of current and previous node"
$ws.Range("G8").Value = "O(N)"
$ws.Range("H8").Value = "O(1)"

# Give the question-link cell the same hyperlink + style treatment as the
# other rows in column C
$ws.Hyperlinks.Add($ws.Range("C8"), "https://leetcode.com/explore/interview/card/top-interview-questions-easy/93/linked-list/560/")
$ws.Range("C8").Style = "Hyperlink"

# Update the view: scroll so column D is left-most and select E10, matching
# where the author ended up after entering the new row
$ws.Range("E10").Select()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
